# Fix some mapping issues in io-model/ItICM (#50)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItICM")

# natural gas and petroleum systems (row 3)
$ws.Range("C3").Value = 0.9
$ws.Range("E3").Value = 0.85
$ws.Range("V3").Value = 0.333

# coal mining (row 6)
$ws.Range("C6").Value = 0.1
$ws.Range("E6").Value = 0.1

# waste management (row 7)
$ws.Range("V7").Value = 0.333

# other industries (row 9)
$ws.Range("E9").Value = 0.05
$ws.Range("V9").Value = 0.334
